$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Price" column (D). These cells store plain text (e.g. "71.115.60"),
# so force text formatting around the write to stop Excel auto-converting the
# string into a number (which would also collapse thousand-separator dots).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "71.115.60"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.867.90"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "697.06"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.15"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.865.94"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.22"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000260"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.519.72"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.864.10"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "71.161.37"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.95"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "497.80"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.95"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.73"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.24"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.17"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.821.41"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.43"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000311"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "163.67"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "49.46"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "417.09"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.61"
$ws.Range("D51").NumberFormat = "General"

# --- Update "Volume(1h)" column (E) ---
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("E3").Value = "  +1.58%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("E7").Value = "  +1.57%  "
$ws.Range("E9").Value = "  +0.10%  "
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("E11").Value = "  -5.60%  "
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("E14").Value = "  +0.91%  "
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("E17").Value = "  +0.47%  "
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("E21").Value = "  -3.11%  "
$ws.Range("E22").Value = "  +3.78%  "
$ws.Range("E24").Value = "  +2.93%  "
$ws.Range("E25").Value = "  +1.38%  "
$ws.Range("E26").Value = "  +2.90%  "
$ws.Range("E27").Value = "  -1.04%  "
$ws.Range("E28").Value = "  -0.55%  "
$ws.Range("E29").Value = "  +2.34%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("E32").Value = "  -0.91%  "
$ws.Range("E33").Value = "  +0.50%  "
$ws.Range("E34").Value = "  +2.63%  "
$ws.Range("E35").Value = "  +0.40%  "
$ws.Range("E36").Value = "  +1.67%  "
$ws.Range("E38").Value = "  +2.01%  "
$ws.Range("E39").Value = "  +9.03%  "
$ws.Range("E40").Value = "  +8.66%  "
$ws.Range("E41").Value = "  -2.88%  "
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  -6.60%  "
$ws.Range("E46").Value = "  +2.22%  "
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("E49").Value = "  +3.71%  "
$ws.Range("E50").Value = "  -1.68%  "
$ws.Range("E51").Value = "  -4.06%  "
